# Update data: 5 June 2021
# Appends the May 2021 (date serial 44317) unemployment observations to both
# the "Canada" sheet (sheet1) and the "Province" sheet (sheet2).

$wb = $excel.ActiveWorkbook
$wsCanada = $wb.Worksheets.Item(1)
$wsProvince = $wb.Worksheets.Item(2)

$newDate = 44317
$dateFormat = "d-mmm-yy"

# ---------------------------------------------------------------------------
# Sheet "Canada": add row 18
# ---------------------------------------------------------------------------
$row = 18

$wsCanada.Cells.Item($row, 1).Value = $newDate
$wsCanada.Cells.Item($row, 1).NumberFormat = $dateFormat

$wsCanada.Cells.Item($row, 2).Value = "Canada"
$wsCanada.Cells.Item($row, 2).NumberFormat = $dateFormat

$wsCanada.Cells.Item($row, 4).Value = 1652.3
$wsCanada.Cells.Item($row, 5).Value = 1093.5999999999999

$wsCanada.Cells.Item($row, 3).Formula = "=(D" + $row + "-E" + $row + ")/E" + $row + "*100"

# ---------------------------------------------------------------------------
# Sheet "Province": add rows 162-171 (Newfoundland & Labrador ... British Columbia)
# ---------------------------------------------------------------------------
$provinceRows = @(
    @{ Row = 162; Name = "Newfoundland & Labrador"; D = 33.700000000000003;  E = 32.9;                  HeaderStyle = $true  },
    @{ Row = 163; Name = "Prince Edward Island";     D = 8.4;                E = 7.6;                   HeaderStyle = $false },
    @{ Row = 164; Name = "Nova Scotia";               D = 48.4;               E = 33;                    HeaderStyle = $false },
    @{ Row = 165; Name = "New Brunswick";             D = 35.799999999999997; E = 29.7;                  HeaderStyle = $false },
    @{ Row = 166; Name = "Quebec";                    D = 296;                E = 225.8;                 HeaderStyle = $false },
    @{ Row = 167; Name = "Ontario";                   D = 733;                E = 410.8;                 HeaderStyle = $false },
    @{ Row = 168; Name = "Manitoba";                  D = 50.8;               E = 35.4;                  HeaderStyle = $false },
    @{ Row = 169; Name = "Saskatchewan";              D = 37.9;               E = 32.299999999999997;    HeaderStyle = $false },
    @{ Row = 170; Name = "Alberta";                   D = 211.7;              E = 164.4;                 HeaderStyle = $false },
    @{ Row = 171; Name = "British Columbia";          D = 196.7;              E = 121.7;                 HeaderStyle = $false }
)

foreach ($item in $provinceRows) {
    $r = $item.Row

    $wsProvince.Cells.Item($r, 1).Value = $newDate
    $wsProvince.Cells.Item($r, 1).NumberFormat = $dateFormat

    $wsProvince.Cells.Item($r, 2).Value = $item.Name
    if ($item.HeaderStyle) {
        $wsProvince.Cells.Item($r, 2).NumberFormat = $dateFormat
    }

    $wsProvince.Cells.Item($r, 4).Value = $item.D
    $wsProvince.Cells.Item($r, 5).Value = $item.E

    $wsProvince.Cells.Item($r, 3).Formula = "=(D" + $r + "-E" + $r + ")/E" + $r + "*100"
}

# ---------------------------------------------------------------------------
# View state: selection / active cell, matching the author's final cursor
# position after data entry (Canada sheet is not the active tab).
# ---------------------------------------------------------------------------
$wsCanada.Range("D19").Select()
$wsProvince.Range("D172").Select()

$wb.Save()
